{"js": "// Fix spelling/wording errors in the second \"-Una vez...\" paragraph:\n//  \"-Una vez dentro veremos unos casilleros con el nombre de las materias\n//   a la nos inscribimos y elegiremos uno par chequear los horarios de\n//   cursada dentro de la secci\u00f3n de inicio \"\n// becomes\n//  \"-Luego de ingresar, dentro veremos unos casilleros con el nombre de\n//   las materias a la que nos inscribimos y entraremos para chequear los\n//   horarios de cursada dentro de la secci\u00f3n de inicio \"\n\nconst body = context.document.body;\n\n// \"-Una vez dentro \" is unique in the document (the first paragraph has\n// \"-Una vez anotados\"), so this only touches the target paragraph.\nconst firstHits = body.search(\"-Una vez dentro \", { matchCase: true });\nfirstHits.load(\"items\");\nawait context.sync();\n\nfor (const hit of firstHits.items) {\n  hit.insertText(\"-Luego de ingresar, dentro \", \"Replace\");\n}\nawait context.sync();\n\n// The second fix, further along in the same paragraph.\nconst secondHits = body.search(\n  \"a la nos inscribimos y elegiremos uno par chequear\",\n  { matchCase: true }\n);\nsecondHits.load(\"items\");\nawait context.sync();\n\nfor (const hit of secondHits.items) {\n  hit.insertText(\"a la que nos inscribimos y entraremos para chequear\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Fix spelling/wording errors in the second \"-Una vez...\" paragraph:\n#  \"-Una vez dentro veremos unos casilleros con el nombre de las materias\n#   a la nos inscribimos y elegiremos uno par chequear los horarios de\n#   cursada dentro de la secci\u00f3n de inicio \"\n# becomes\n#  \"-Luego de ingresar, dentro veremos unos casilleros con el nombre de\n#   las materias a la que nos inscribimos y entraremos para chequear los\n#   horarios de cursada dentro de la secci\u00f3n de inicio \"\n\n$d = $word.ActiveDocument\n\n# \"-Una vez dentro \" is unique in the document (the first paragraph reads\n# \"-Una vez anotados\"), so this Find/Replace only touches the target spot.\n$find1 = $d.Content.Find\n$find1.Text = \"-Una vez dentro \"\n$find1.Replacement.Text = \"-Luego de ingresar, dentro \"\n$find1.Execute(\n    $find1.Text,       # FindText\n    $false,            # MatchCase\n    $false,            # MatchWholeWord\n    $false,            # MatchWildcards\n    $false,            # MatchSoundsLike\n    $false,            # MatchAllWordForms\n    $true,             # Forward\n    1,                 # Wrap (wdFindContinue)\n    $false,            # Format\n    $find1.Replacement.Text, # ReplaceWith\n    2                  # Replace (wdReplaceAll)\n)\n\n# Second fix, later in the same paragraph.\n$find2 = $d.Content.Find\n$find2.Text = \"a la nos inscribimos y elegiremos uno par chequear\"\n$find2.Replacement.Text = \"a la que nos inscribimos y entraremos para chequear\"\n$find2.Execute(\n    $find2.Text,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $find2.Replacement.Text,\n    2\n)\n"}
